$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) updates to column F (想去人数)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3422
$ws1.Range("F5").Value = 6975
$ws1.Range("F6").Value = 2451
$ws1.Range("F7").Value = 43
$ws1.Range("F8").Value = 112
$ws1.Range("F14").Value = 573

# Sheet "全部类型" (All Types) updates to column F (想去人数)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3422
$ws4.Range("F6").Value = 6975
$ws4.Range("F7").Value = 2451
$ws4.Range("F8").Value = 43
$ws4.Range("F9").Value = 112
$ws4.Range("F15").Value = 573
